$wb = $excel.ActiveWorkbook

# Add a new worksheet named VERSION, positioned after "Staff Data" and before "Sheet2"
$staffSheet = $wb.Worksheets.Item("Staff Data")
$sheet2 = $wb.Worksheets.Item("Sheet2")
$versionSheet = $wb.Worksheets.Add($sheet2)
$versionSheet.Name = "VERSION"

$versionSheet.Range("A6").Value = "File version"
$versionSheet.Range("B6").Value = "vx.xx"

$versionSheet.Select() | Out-Null
$versionSheet.Range("B6").Select() | Out-Null
